$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 0.208838632887556
$ws.Range("C3").Value = 0.178857494432493
$ws.Range("C4").Value = 0.187422136325044
$ws.Range("C5").Value = 0.188547853061829
$ws.Range("C6").Value = 0.236333883293078
